$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C2:C7) from 2023-10-25 (serial 45224) to 2023-11-03 (serial 45233)
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 45224) {
        $cell.Value = 45233
    }
}
